# Assignment_0.docx edit:
# Replace the old "Finally, go ahead..." sentence that introduced
# Assignment 1 with new instructions about running/submitting 0.ss,
# while leaving the existing "PLC Grading Server" hyperlink + trailing
# period untouched.

$d = $word.ActiveDocument

$oldText = "Finally, go ahead and get started on Assignment 1 and practice submitting your code via the "

# Locate the exact run text via Find so we don't have to hard-code offsets.
$target = $d.Content
$found = $target.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the original sentence to replace."
}

# Remove the old text first, leaving a collapsed insertion point in its place.
$target.Delete()

# Insert the three replacement runs. Inserting repeatedly *before* the same
# (still-collapsed) range, in reverse order, keeps them as separate <w:r>
# runs (rather than being coalesced into a single run) and each new run
# picks up the plain-paragraph formatting instead of the neighbouring
# hyperlink's formatting.
$target.InsertBefore(" ")
$target.InsertBefore(", run in your scheme interpter to verify it works, and them submit it to the")
$target.InsertBefore("Take the already written code in 0.ss")
